$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 336
$ws.Range("L2").Value = "stimuli/img_njhlh.png"
$ws.Range("M2").Value = 59.74418604651163
$ws.Range("N2").Value = 41.51162790697674
$ws.Range("O2").Value = 50.62790697674419
$ws.Range("P2").Value = 43
$ws.Range("Q2").Value = 4
$ws.Range("R2").Value = 4
$ws.Range("S2").Value = 4

$ws.Range("F3").Value = 337
$ws.Range("L3").Value = "stimuli/img_vnxft.png"
$ws.Range("M3").Value = 53.22727272727273
$ws.Range("N3").Value = 34.84090909090909
$ws.Range("O3").Value = 44.03409090909091
$ws.Range("P3").Value = 44
$ws.Range("Q3").Value = 3
$ws.Range("R3").Value = 3
$ws.Range("S3").Value = 3

$ws.Range("F4").Value = 338
$ws.Range("L4").Value = "stimuli/img_j4ttn.png"
$ws.Range("M4").Value = 12.61904761904762
$ws.Range("N4").Value = 11.42857142857143
$ws.Range("O4").Value = 12.02380952380952
$ws.Range("P4").Value = 42
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 1
$ws.Range("S4").Value = 1

$ws.Range("F5").Value = 339

$ws.Range("F6").Value = 340
$ws.Range("L6").Value = "stimuli/img_tujn3.png"
$ws.Range("M6").Value = 81.40909090909091
$ws.Range("N6").Value = 62.52272727272727
$ws.Range("O6").Value = 71.96590909090909
$ws.Range("P6").Value = 44
$ws.Range("Q6").Value = 8
$ws.Range("R6").Value = 8
$ws.Range("S6").Value = 8

$ws.Range("F7").Value = 341
$ws.Range("H7").Value = "bedrooms"
$ws.Range("I7").Value = "distractor"
$ws.Range("K7").Value = "f"
$ws.Range("L7").Value = "stimuli/img_mdpr4.png"
$ws.Range("M7").Value = 74.04255319148936
$ws.Range("N7").Value = 54.70212765957447
$ws.Range("O7").Value = 64.37234042553192
$ws.Range("Q7").Value = 6
$ws.Range("R7").Value = 6
$ws.Range("S7").Value = 6

$ws.Range("F8").Value = 342
$ws.Range("L8").Value = "stimuli/img_ra2nm.png"
$ws.Range("M8").Value = 70.75
$ws.Range("N8").Value = 50.375
$ws.Range("O8").Value = 60.5625
$ws.Range("P8").Value = 40
$ws.Range("Q8").Value = 6
$ws.Range("R8").Value = 6
$ws.Range("S8").Value = 6

$ws.Range("F9").Value = 343
$ws.Range("L9").Value = "stimuli/img_vh7v8.png"
$ws.Range("M9").Value = 78.70454545454545
$ws.Range("N9").Value = 59.63636363636363
$ws.Range("O9").Value = 69.17045454545455
$ws.Range("P9").Value = 44
$ws.Range("Q9").Value = 7
$ws.Range("R9").Value = 7
$ws.Range("S9").Value = 7

$ws.Range("F10").Value = 344
$ws.Range("H10").Value = "bedrooms"
$ws.Range("I10").Value = "distractor"
$ws.Range("K10").Value = "f"
$ws.Range("L10").Value = "stimuli/img_e26ut.png"
$ws.Range("M10").Value = 81.07692307692308
$ws.Range("N10").Value = 61.28205128205128
$ws.Range("O10").Value = 71.17948717948718
$ws.Range("P10").Value = 39
$ws.Range("Q10").Value = 8
$ws.Range("R10").Value = 8
$ws.Range("S10").Value = 8

$ws.Range("F11").Value = 345
$ws.Range("L11").Value = "stimuli/img_73pyk.png"
$ws.Range("M11").Value = 69.27659574468085
$ws.Range("N11").Value = 47.27659574468085
$ws.Range("O11").Value = 58.27659574468085
$ws.Range("Q11").Value = 5
$ws.Range("R11").Value = 5
$ws.Range("S11").Value = 5

$ws.Range("F12").Value = 346
$ws.Range("H12").Value = "bedrooms"
$ws.Range("I12").Value = "distractor"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/img_d3t0o.png"
$ws.Range("M12").Value = 66.95121951219512
$ws.Range("N12").Value = 42.92682926829269
$ws.Range("O12").Value = 54.9390243902439
$ws.Range("P12").Value = 41
$ws.Range("Q12").Value = 4
$ws.Range("R12").Value = 4
$ws.Range("S12").Value = 4

$ws.Range("F13").Value = 347
$ws.Range("L13").Value = "stimuli/img_swq34.png"
$ws.Range("M13").Value = 64.11363636363636
$ws.Range("N13").Value = 43.04545454545455
$ws.Range("O13").Value = 53.57954545454545
$ws.Range("P13").Value = 44
$ws.Range("Q13").Value = 5
$ws.Range("R13").Value = 5
$ws.Range("S13").Value = 5

$ws.Range("F14").Value = 348
$ws.Range("H14").Value = "living_rooms"
$ws.Range("I14").Value = "target"
$ws.Range("K14").Value = "j"
$ws.Range("L14").Value = "stimuli/img_g13d5.png"
$ws.Range("M14").Value = 73
$ws.Range("N14").Value = 51.51111111111111
$ws.Range("O14").Value = 62.25555555555556
$ws.Range("P14").Value = 45
$ws.Range("Q14").Value = 6
$ws.Range("R14").Value = 6
$ws.Range("S14").Value = 6

$ws.Range("F15").Value = 349
$ws.Range("L15").Value = "stimuli/img_ac0ey.png"
$ws.Range("M15").Value = 86.62222222222222
$ws.Range("N15").Value = 70.02222222222223
$ws.Range("O15").Value = 78.32222222222222
$ws.Range("P15").Value = 45
$ws.Range("Q15").Value = 9
$ws.Range("R15").Value = 9
$ws.Range("S15").Value = 9

$ws.Range("F16").Value = 350
$ws.Range("H16").Value = "living_rooms"
$ws.Range("I16").Value = "target"
$ws.Range("K16").Value = "j"
$ws.Range("L16").Value = "stimuli/img_of8d6.png"
$ws.Range("M16").Value = 26.04878048780488
$ws.Range("N16").Value = 19.14634146341463
$ws.Range("O16").Value = 22.59756097560975
$ws.Range("P16").Value = 41
$ws.Range("Q16").Value = 1
$ws.Range("R16").Value = 1
$ws.Range("S16").Value = 1

$ws.Range("F17").Value = 351
$ws.Range("L17").Value = "stimuli/img_rg4in.png"
$ws.Range("M17").Value = 49.3695652173913
$ws.Range("N17").Value = 30.21739130434782
$ws.Range("O17").Value = 39.79347826086956
$ws.Range("P17").Value = 46
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 3
$ws.Range("S17").Value = 3

$ws.Range("F18").Value = 352
$ws.Range("L18").Value = "stimuli/img_jp28n.png"
$ws.Range("M18").Value = 65.02564102564102
$ws.Range("N18").Value = 44.97435897435897
$ws.Range("O18").Value = 55
$ws.Range("P18").Value = 39
$ws.Range("Q18").Value = 4
$ws.Range("R18").Value = 4
$ws.Range("S18").Value = 4

$ws.Range("F19").Value = 353
$ws.Range("L19").Value = "stimuli/img_9oofc.png"
$ws.Range("M19").Value = 82.47619047619048
$ws.Range("N19").Value = 65.5
$ws.Range("O19").Value = 73.98809523809524
$ws.Range("P19").Value = 42
$ws.Range("Q19").Value = 8
$ws.Range("R19").Value = 8
$ws.Range("S19").Value = 8

$ws.Range("F20").Value = 354
$ws.Range("L20").Value = "stimuli/img_rru0v.png"
$ws.Range("M20").Value = 56.45238095238095
$ws.Range("N20").Value = 39.42857142857143
$ws.Range("O20").Value = 47.94047619047619
$ws.Range("P20").Value = 42
$ws.Range("Q20").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("S20").Value = 4

$ws.Range("F21").Value = 355
$ws.Range("H21").Value = "living_rooms"
$ws.Range("I21").Value = "target"
$ws.Range("K21").Value = "j"
$ws.Range("L21").Value = "stimuli/img_vgh2g.png"
$ws.Range("M21").Value = 93.81395348837209
$ws.Range("N21").Value = 78.27906976744185
$ws.Range("O21").Value = 86.04651162790697
$ws.Range("P21").Value = 43
$ws.Range("Q21").Value = 10
$ws.Range("R21").Value = 10
$ws.Range("S21").Value = 10

$ws.Range("F22").Value = 356
$ws.Range("H22").Value = "living_rooms"
$ws.Range("I22").Value = "target"
$ws.Range("K22").Value = "j"
$ws.Range("L22").Value = "stimuli/img_5jp4f.png"
$ws.Range("M22").Value = 84.85714285714286
$ws.Range("N22").Value = 67.83333333333333
$ws.Range("O22").Value = 76.3452380952381
$ws.Range("P22").Value = 42
$ws.Range("Q22").Value = 9
$ws.Range("R22").Value = 9
$ws.Range("S22").Value = 9

$ws.Range("F23").Value = 357
$ws.Range("L23").Value = "stimuli/img_xr3up.png"
$ws.Range("M23").Value = 76.24444444444444
$ws.Range("N23").Value = 55.88888888888889
$ws.Range("O23").Value = 66.06666666666666
$ws.Range("P23").Value = 45
$ws.Range("Q23").Value = 7
$ws.Range("R23").Value = 7
$ws.Range("S23").Value = 7

$ws.Range("F24").Value = 358
$ws.Range("L24").Value = "stimuli/img_rych7.png"
$ws.Range("M24").Value = 30.4468085106383
$ws.Range("N24").Value = 23.4468085106383
$ws.Range("O24").Value = 26.9468085106383
$ws.Range("P24").Value = 47
$ws.Range("Q24").Value = 2
$ws.Range("R24").Value = 2
$ws.Range("S24").Value = 2

$ws.Range("F25").Value = 359
$ws.Range("H25").Value = "bedrooms"
$ws.Range("I25").Value = "distractor"
$ws.Range("K25").Value = "f"
$ws.Range("L25").Value = "stimuli/img_91csq.png"
$ws.Range("M25").Value = 50.44736842105263
$ws.Range("N25").Value = 28.34210526315789
$ws.Range("O25").Value = 39.39473684210526
$ws.Range("P25").Value = 38
$ws.Range("Q25").Value = 2
$ws.Range("R25").Value = 2
$ws.Range("S25").Value = 2

$ws.Range("F26").Value = 360
$ws.Range("L26").Value = "stimuli/img_gxm46.png"
$ws.Range("M26").Value = 74.78378378378379
$ws.Range("N26").Value = 54
$ws.Range("O26").Value = 64.3918918918919
$ws.Range("P26").Value = 37
$ws.Range("Q26").Value = 6
$ws.Range("R26").Value = 6
$ws.Range("S26").Value = 6

$ws.Range("F27").Value = 361
$ws.Range("L27").Value = "stimuli/img_syam3.png"
$ws.Range("M27").Value = 41.32432432432432
$ws.Range("N27").Value = 26.2972972972973
$ws.Range("O27").Value = 33.81081081081081
$ws.Range("P27").Value = 37
$ws.Range("Q27").Value = 2
$ws.Range("R27").Value = 2
$ws.Range("S27").Value = 2
